$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "Golang Architect / Principal Backend Architect - Atlanta, GA -onsite"
$ws.Range("B30").Value = "https://www.dice.com/job-detail/8c7df8da-e8d3-48e0-b8a4-fa0e0d0fd9c9"
$ws.Range("C30").Value = "Atlanta, Georgia"
$ws.Range("D30").Value = "Contract"
$ws.Range("E30").Value = "Depends on Experience"
$ws.Range("F30").Value = "Oxford Global Resources"
